$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
Write-Output "before: $($sh.TextFrame.TextRange.Text)"
try {
  $sh.TextFrame.TextRange.InsertDateTime(2, 0)
  Write-Output "after: $($sh.TextFrame.TextRange.Text)"
} catch {
  Write-Output "ERR: $_"
}
